# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" locale sheets:
#  - Status column (B) flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" for the two real source-file rows.
#  - Two new columns are populated: "Latest Target File" (E) and
#    "Latest Handback File" (F), mirroring the source .md file (A) and the
#    handoff .xlf file (C) respectively, each as a hyperlinked file name.
#  - "Latest Handback DateTime" (G) moves from the placeholder
#    0001-01-01 00:00:00 to the real handback timestamp.

$wb = $excel.ActiveWorkbook

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackTime2,
        [string]$HandbackTime3
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Build a (row,col) -> target-url lookup from the existing hyperlinks,
    # since Hyperlinks.Item(n) isn't reliable on this host.
    $urlMap = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $key = [string]$hl.Range.Row + "_" + [string]$hl.Range.Column
        $urlMap[$key] = $hl.Address
    }

    $mdName1  = $ws.Range("A2").Value2
    $xlfName1 = $ws.Range("C2").Value2
    $mdName2  = $ws.Range("A3").Value2
    $xlfName2 = $ws.Range("C3").Value2

    $mdUrl1  = $urlMap["2_1"]
    $xlfUrl1 = $urlMap["2_3"]
    $mdUrl2  = $urlMap["3_1"]
    $xlfUrl2 = $urlMap["3_3"]

    # Status -> handed back
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Latest Target File (E) / Latest Handback File (F) for row 2
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", $mdName1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl1, "", "", $xlfName1) | Out-Null
    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"

    # Latest Target File (E) / Latest Handback File (F) for row 3
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", $mdName2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl2, "", "", $xlfName2) | Out-Null
    $ws.Range("E3").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"

    # Latest Handback DateTime (G)
    $ws.Range("G2").Value = $HandbackTime2
    $ws.Range("G3").Value = $HandbackTime3
}

Update-LocaleSheet "zh-cn" "2016-02-18 04:16:04" "2016-02-18 04:16:04"
Update-LocaleSheet "de-de" "2016-02-18 04:16:27" "2016-02-18 04:16:27"

# The "Ready for handoff" status also shows on the Overview roll-up sheet
# (same shared-string text, now retired in favor of the handback wording).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"
